$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name (column B) updates ---
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("B10").Value = 'Polygon'
$ws.Range("B11").Value = 'Solana'
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("B15").Value = 'TRON'
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("B20").Value = 'Dai'
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("B26").Value = 'Monero'
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("B32").Value = 'Stellar'
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("B36").Value = 'Hedera'
$ws.Range("B37").Value = 'VeChain'
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("B39").Value = 'Frax'
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("B42").Value = 'Aptos'
$ws.Range("B43").Value = 'Algorand'
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("B48").Value = 'Cronos'
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("B50").Value = 'Quant'
$ws.Range("B51").Value = 'Elrond'

# --- Link (column C) updates ---
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C39").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'

# --- Price (column D) updates that are unambiguous text (multiple '.' separators) ---
$ws.Range("D2").Value = '28.666.29'
$ws.Range("D3").Value = '1.851.49'
$ws.Range("D12").Value = '1.912.37'
$ws.Range("D21").Value = '28.684.30'
$ws.Range("D25").Value = '2.111.57'

# --- Price (column D) updates that Excel would otherwise parse as numbers; force text ---
$ambiguousPriceCells = @(
    @("D4", '1.002'),
    @("D5", '335.51'),
    @("D7", '0.4661'),
    @("D8", '0.3911'),
    @("D9", '0.07914'),
    @("D10", '0.9857'),
    @("D11", '22.32'),
    @("D13", '5.852'),
    @("D14", '7.015'),
    @("D15", '0.06849'),
    @("D16", '87.78'),
    @("D17", '1.003'),
    @("D18", '0.00001007'),
    @("D19", '17.16'),
    @("D20", '1.003'),
    @("D22", '5.405'),
    @("D23", '11.30'),
    @("D24", '2.140'),
    @("D26", '153.21'),
    @("D27", '19.48'),
    @("D28", '6.051'),
    @("D29", '2.029'),
    @("D30", '117.78'),
    @("D31", '0.9775'),
    @("D32", '0.09432'),
    @("D33", '5.379'),
    @("D34", '3.484'),
    @("D35", '1.353'),
    @("D36", '0.06172'),
    @("D37", '0.02201'),
    @("D38", '1.163'),
    @("D39", '1.002'),
    @("D40", '0.5733'),
    @("D41", '7.649'),
    @("D42", '10.23'),
    @("D43", '0.1804'),
    @("D44", '2.384'),
    @("D45", '1.247'),
    @("D46", '0.5400'),
    @("D47", '11.74'),
    @("D48", '0.07143'),
    @("D49", '1.909'),
    @("D50", '114.94'),
    @("D51", '43.54')
)
foreach ($pair in $ambiguousPriceCells) {
    $cellRef = $pair[0]
    $cellVal = $pair[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $cellVal
    $rng.Style = "Normal"
}

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = '  -3.25%  '
$ws.Range("E3").Value = '  -3.99%  '
$ws.Range("E4").Value = '  -1.13%  '
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("E7").Value = '  -3.38%  '
$ws.Range("E8").Value = '  -3.82%  '
$ws.Range("E9").Value = '  -4.03%  '
$ws.Range("E10").Value = '  -2.74%  '
$ws.Range("E11").Value = '  -6.84%  '
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("E13").Value = '  -4.16%  '
$ws.Range("E14").Value = '  -3.58%  '
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("E16").Value = '  -4.52%  '
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("E18").Value = '  -3.21%  '
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("E21").Value = '  -3.12%  '
$ws.Range("E22").Value = '  -4.96%  '
$ws.Range("E23").Value = '  -5.85%  '
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("E26").Value = '  -1.97%  '
$ws.Range("E27").Value = '  -2.97%  '
$ws.Range("E28").Value = '  -6.36%  '
$ws.Range("E29").Value = '  -3.28%  '
$ws.Range("E30").Value = '  -2.60%  '
$ws.Range("E31").Value = '  -3.77%  '
$ws.Range("E32").Value = '  -2.23%  '
$ws.Range("E33").Value = '  -4.58%  '
$ws.Range("E34").Value = '  -2.64%  '
$ws.Range("E35").Value = '  -2.05%  '
$ws.Range("E36").Value = '  -3.67%  '
$ws.Range("E37").Value = '  -4.38%  '
$ws.Range("E38").Value = '  -2.26%  '
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("E40").Value = '  -3.88%  '
$ws.Range("E41").Value = '  -3.12%  '
$ws.Range("E42").Value = '  -5.13%  '
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("E44").Value = '  -2.68%  '
$ws.Range("E45").Value = '  -2.99%  '
$ws.Range("E46").Value = '  -3.14%  '
$ws.Range("E47").Value = '  -5.31%  '
$ws.Range("E48").Value = '  -5.47%  '
$ws.Range("E49").Value = '  -2.64%  '
$ws.Range("E50").Value = '  -3.87%  '
$ws.Range("E51").Value = '  +3.24%  '
